# Insert a new "2022-Q3" sheet before "2022-Q2", fill it with fund-holding
# data, and prepend a matching summary row on "总计".

$wb = $excel.ActiveWorkbook

$oldQ2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($oldQ2)
$newSheet.Name = "2022-Q3"

$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @("900090", "中信卓越成长两年持有期混合B", "51.50", "93.14", "4.03", "2.0754", 6),
    @("900010", "中信卓越成长两年持有期混合A", "14.24", "93.14", "4.03", "0.5739", 6),
    @("900100", "中信卓越成长两年持有期混合C", "4.61", "93.14", "4.03", "0.1858", 6),
    @("470888", "汇添富香港优势精选混合（QDII）", "1.63", "78.50", "3.45", "0.0562", 9)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2
    $newSheet.Cells.Item($excelRow, 1).Value = $r
    for ($c = 0; $c -lt $row.Length; $c++) {
        $newSheet.Cells.Item($excelRow, $c + 2).Value = $row[$c]
    }
}

# Prepend a "2022-Q3" row on the summary sheet, pushing existing rows down.
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:D2").Insert()
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 2.89
